$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "81÷8=10, 1"
$t.Cell(1,2).Range.Text = "51÷4=12, 3"
$t.Cell(1,3).Range.Text = "78÷2=39, 0"
$t.Cell(1,4).Range.Text = "24÷7=3, 3"
$t.Cell(1,5).Range.Text = "39÷5=7, 4"
$t.Cell(5,1).Range.Text = "57÷8=7, 1"
$t.Cell(5,2).Range.Text = "43÷8=5, 3"
$t.Cell(5,3).Range.Text = "58÷5=11, 3"
$t.Cell(5,4).Range.Text = "35÷5=7, 0"
$t.Cell(5,5).Range.Text = "71÷4=17, 3"
$t.Cell(9,1).Range.Text = "80÷4=20, 0"
$t.Cell(9,2).Range.Text = "51÷7=7, 2"
$t.Cell(9,3).Range.Text = "81÷4=20, 1"
$t.Cell(9,4).Range.Text = "83÷8=10, 3"
$t.Cell(9,5).Range.Text = "79÷6=13, 1"
$t.Cell(13,1).Range.Text = "39÷9=4, 3"
$t.Cell(13,2).Range.Text = "69÷6=11, 3"
$t.Cell(13,3).Range.Text = "88÷8=11, 0"
$t.Cell(13,4).Range.Text = "67÷7=9, 4"
$t.Cell(13,5).Range.Text = "36÷2=18, 0"
$t.Cell(17,1).Range.Text = "39÷6=6, 3"
$t.Cell(17,2).Range.Text = "29÷5=5, 4"
$t.Cell(17,3).Range.Text = "65÷6=10, 5"
$t.Cell(17,4).Range.Text = "76÷4=19, 0"
$t.Cell(17,5).Range.Text = "15÷5=3, 0"
